$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 ("cargo") was entered as plain/unstyled text; normalize it
#     to match the rest of the formatted table (numbers + same styles). ---
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "cargo"
$ws.Range("C17").Value = 500
$ws.Range("D17").Value = 10

$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("B17").Style = "Normal"

# --- Row 18 ("Saree") new row, styled ---
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Saree"
$ws.Range("C18").Value = 1234
$ws.Range("D18").Value = 5
$ws.Rows.Item(18).RowHeight = 18.75

# --- Row 19 ("saree") new row, plain/unstyled ---
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "saree"
$ws.Range("C19").Value = "1111"
$ws.Range("D19").Value = "222"

# --- Row 20 ("Shirt") new row, plain/unstyled ---
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Shirt"
$ws.Range("C20").Value = "120"
$ws.Range("D20").Value = "34"
